# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx price/volume/coin-name refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.401.15"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "1.722.97"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4873"
$ws.Range("D7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2586"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06191"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.47%  "

$ws.Range("D10").Value = "1.730.85"
$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06964"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.46"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.540"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5970"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.22"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "26.400.37"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9999"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007219"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.32"
$ws.Range("D20").ClearFormats()

$ws.Range("D21").Value = "1.944.60"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.443"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.484"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.105"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.69%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.24"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.401"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.68"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.724"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.926"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08008"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.668"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04492"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.603"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9976"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6246"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9376"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.41%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.949"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.55%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.385"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01473"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.05"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.325"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3834"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.856"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1163"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.75%  "

$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.12"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.36%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.704"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.228"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.82"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.27%  "

Write-Host "Applied 100 cell updates"
